# Minutes of 28th meeting: rework the "3. Implementation" subtasks (rows
# 17-26) into 16 rows (17-32) with new activity names/dates, which pushes
# everything below (old rows 18-38) down by 6 rows (now 24-44).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Planner")

# 1. Make room: insert 6 new rows right after row 17 (old rows 18-26 -> 24-32,
#    old rows 27-38 -> 33-44). Excel copies formatting from the row above.
$ws.Rows("18:23").Insert()

# 2. Row 17 keeps its activity name but now has actual-start/duration data.
$ws.Range("B17").Value = "3.1 Create testing structure"
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 14
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 1

# 3. New / rewritten implementation subtasks (rows 18-28).
$ws.Range("B18").Value = "3.2 Learn related knowledge"
$ws.Range("C18").Value = 15
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 15
$ws.Range("F18").Value = 2
$ws.Range("G18").Value = 1

$ws.Range("B19").Value = "3.3 Make basic components"
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 16
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 1

$ws.Range("B20").Value = "3.4 UI for start page & subpages"
$ws.Range("C20").Value = 16
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 16
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 1

$ws.Range("B21").Value = "3.5 make swap animation"
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 16
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 1

$ws.Range("B22").Value = "3.6 Bubble sort algorithm"
$ws.Range("C22").Value = 17
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 17
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 1

$ws.Range("B23").Value = "3.7 Router feat & local storage"
$ws.Range("C23").Value = 17
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 17
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 1

$ws.Range("B24").Value = "3.8 Animate 2 algorithms (SI)"
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = 17
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = 1

$ws.Range("B25").Value = "3.9 Help and setting"
$ws.Range("C25").Value = 17
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 17
$ws.Range("F25").Value = $null
$ws.Range("G25").Value = 0.8

$ws.Range("B26").Value = "3.10 Correctness proposal"
$ws.Range("C26").Value = 18
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 18
$ws.Range("F26").Value = $null
$ws.Range("G26").Value = 0.8

$ws.Range("B27").Value = "3.11 improve ui with electron"
$ws.Range("C27").Value = 18
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 18
$ws.Range("F27").Value = $null
$ws.Range("G27").Value = $null

$ws.Range("B28").Value = "3.12 animate merge and heap"
$ws.Range("C28").Value = 18
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 18
$ws.Range("F28").Value = $null
$ws.Range("G28").Value = $null

# 4. Update the conditional-formatting ranges that previously ended at the old
#    last data row (38) / total row (39) so they now cover the new rows.
$cfs = $ws.Cells.FormatConditions
for ($i = 1; $i -le $cfs.Count(); $i++) {
    $cf = $cfs.Item($i)
    $addr = $cf.AppliesTo().Address()
    if ($addr -eq "`$H`$5:`$AI`$38") {
        $cf.ModifyAppliesToRange($ws.Range("H5:AI44"))
    } elseif ($addr -eq "`$B`$39:`$BO`$39") {
        $cf.ModifyAppliesToRange($ws.Range("B45:BO45"))
    }
}

# 5. Leave the view where the author last left it.
$ws.Activate() | Out-Null
$ws.Range("F25").Select() | Out-Null
